$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(318294931, Shalev  Afanasenko: 3,-9)"
$ws.Range("B1").Value = "(305487936, Avihai  Kipnis: 8,-1)"
$ws.Range("C1").Value = "(313227928, Aviv  Levi: 9,-1)"
$ws.Range("D1").Value = "(205807308, Sariel  Basis: 2,8)"
$ws.Range("E1").Value = "(315891549, Raz  Halaby: 3,-7)"
$ws.Range("F1").Value = "(315060103, Dan  Mshelh: -3,-6)"
$ws.Range("G1").Value = "(313925141, Elad   Amer: 0,-1)"

$ws.Range("A3").Value = "cost: 453.81152608925026"
$ws.Range("A4").Value = "time: 61.25878944132147"
